# PQ_Challenge_196 - "Looking an alternative excel solution"
# Adds two new sheets (EDA, Alt) that are copies of the original
# "Original" sheet, and adds an alternative LET/LAMBDA based array
# formula solution to the "Alt" sheet.

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item("Original")

# Duplicate "Original" -> "EDA" placed right after "Original"
$orig.Copy($null, $orig)
$eda = $wb.Worksheets.Item($orig.Index + 1)
$eda.Name = "EDA"

# Duplicate "Original" again -> "Alt" placed right after "EDA"
$orig.Copy($null, $eda)
$alt = $wb.Worksheets.Item($eda.Index + 1)
$alt.Name = "Alt"

# Add the alternative LET/LAMBDA array-formula solution on the Alt sheet
# (Transposes Class/Subject/Marks into a Class-by-Subject grid, the same
# result produced by the Original/EDA "helper columns" approach.)
$alt.Activate()
$altFormula = '=LET(a,A2:A11,b,B2:B11,u,UNIQUE(a),s,TOROW(SORT(UNIQUE(b))),e,LAMBDA(x,MAKEARRAY(ROWS(u),COLUMNS(s),LAMBDA(r,c,IF(SUMPRODUCT((a=INDEX(u,r,1))*(b=INDEX(s,1,c)))=0,"",SUMPRODUCT((a=INDEX(u,r,1))*(b=INDEX(s,1,c))*x))))),HSTACK(e(a),e(C2:C11)))'
$alt.Range("B17").Formula2 = $altFormula

$alt.Range("J8").Select()
